$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values were recalculated (K = something derived instead of old "Strike#")
$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
